# Auto-generated Excel COM-interop script
# Applies the value corrections described in the commit diff
# across all 8 profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 820.6
$ws.Range("I38").Value = 820.6
$ws.Range("K38").Value = 2461.8
$ws.Range("M38").Value = -2089.8
$ws.Range("H58").Value = 427.25
$ws.Range("I58").Value = 427.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1281.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1131.75
$ws.Range("N58").ClearContents()
$ws.Range("H137").Value = 2538.9167
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2538.9167
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 7616.750100000001
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -12716.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1346.3954
$ws.Range("I2").Value = 1028.3529
$ws.Range("K2").Value = 1028.3529
$ws.Range("M2").Value = -915.3529000000001
$ws.Range("H35").Value = 9072.143
$ws.Range("J35").Value = 27200
$ws.Range("L35").Value = 27200
$ws.Range("N35").Value = -28012
$ws.Range("H63").Value = 2793.3
$ws.Range("I63").Value = 3133.2856
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 3133.2856
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -2447.2856
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2793.3
$ws.Range("I66").Value = 3133.2856
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 15666.428
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -12234.428
$ws.Range("N66").Value = -16864
$ws.Range("H116").Value = 1346.3954
$ws.Range("I116").Value = 1028.3529
$ws.Range("K116").Value = 1028.3529
$ws.Range("M116").Value = 1265.6471
$ws.Range("H132").Value = 2173.543
$ws.Range("I132").Value = 1293.7222
$ws.Range("K132").Value = 3881.1666
$ws.Range("M132").Value = -1351.1666
$ws.Range("H134").Value = 79999.664
$ws.Range("J134").Value = 79999.664
$ws.Range("L134").Value = 79999.664
$ws.Range("N134").Value = -90139.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1346.3954
$ws.Range("I3").Value = 1028.3529
$ws.Range("K3").Value = 1028.3529
$ws.Range("M3").Value = -914.3529000000001
$ws.Range("H8").Value = 8857.571
$ws.Range("I8").Value = 8751
$ws.Range("J8").Value = 8999.666999999999
$ws.Range("K8").Value = 8751
$ws.Range("L8").Value = 8999.666999999999
$ws.Range("M8").Value = -8611
$ws.Range("N8").Value = -9279.666999999999
$ws.Range("H54").Value = 1500
$ws.Range("I54").Value = 1500
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1500
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1016
$ws.Range("N54").ClearContents()
$ws.Range("H134").Value = 4560.4287
$ws.Range("I134").Value = 4560.4287
$ws.Range("K134").Value = 13681.2861
$ws.Range("M134").Value = -11146.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2327835.2
$ws.Range("I31").Value = 3032486.8
$ws.Range("J31").Value = 2485
$ws.Range("K31").Value = 3032486.8
$ws.Range("L31").Value = 2485
$ws.Range("M31").Value = -3032191.8
$ws.Range("N31").Value = -3075
$ws.Range("H34").Value = 2327835.2
$ws.Range("I34").Value = 3032486.8
$ws.Range("J34").Value = 2485
$ws.Range("K34").Value = 3032486.8
$ws.Range("L34").Value = 2485
$ws.Range("M34").Value = -3032284.8
$ws.Range("N34").Value = -2889
$ws.Range("H59").Value = 36332.668
$ws.Range("I59").Value = 24999
$ws.Range("J59").Value = 41999.5
$ws.Range("K59").Value = 24999
$ws.Range("L59").Value = 41999.5
$ws.Range("M59").Value = -23854
$ws.Range("N59").Value = -44289.5
$ws.Range("H60").Value = 35039.6
$ws.Range("I60").Value = 36899
$ws.Range("J60").Value = 34833
$ws.Range("K60").Value = 36899
$ws.Range("L60").Value = 34833
$ws.Range("M60").Value = -36388
$ws.Range("N60").Value = -35855
$ws.Range("H107").Value = 660.8823
$ws.Range("I107").Value = 540.5
$ws.Range("K107").Value = 540.5
$ws.Range("M107").Value = 1379.5
$ws.Range("H138").Value = 109441.9
$ws.Range("J138").Value = 109441.9
$ws.Range("L138").Value = 109441.9
$ws.Range("N138").Value = -119721.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1699.75
$ws.Range("I132").Value = 1856.5
$ws.Range("J132").Value = 1587.7858
$ws.Range("K132").Value = 16708.5
$ws.Range("L132").Value = 14290.0722
$ws.Range("M132").Value = -14178.5
$ws.Range("N132").Value = -19350.0722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 6022006
$ws.Range("I44").Value = 31000
$ws.Range("K44").Value = 31000
$ws.Range("M44").Value = -30404
$ws.Range("H107").Value = 269.66666
$ws.Range("I107").Value = 164.4
$ws.Range("J107").Value = 401.25
$ws.Range("K107").Value = 164.4
$ws.Range("L107").Value = 401.25
$ws.Range("M107").Value = 1755.6
$ws.Range("N107").Value = -4241.25
$ws.Range("H113").Value = 1579.2
$ws.Range("I113").Value = 1511
$ws.Range("J113").Value = 1596.25
$ws.Range("K113").Value = 1511
$ws.Range("L113").Value = 1596.25
$ws.Range("M113").Value = 659
$ws.Range("N113").Value = -5936.25
$ws.Range("H132").Value = 9649.643
$ws.Range("I132").Value = 10238.077
$ws.Range("K132").Value = 30714.231
$ws.Range("M132").Value = -28184.231
$ws.Range("H139").Value = 104999.336
$ws.Range("J139").Value = 104999.336
$ws.Range("L139").Value = 104999.336
$ws.Range("N139").Value = -115279.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1148.2609
$ws.Range("I55").Value = 485.625
$ws.Range("J55").Value = 1501.6666
$ws.Range("K55").Value = 485.625
$ws.Range("L55").Value = 1501.6666
$ws.Range("M55").Value = -312.625
$ws.Range("N55").Value = -1847.6666
$ws.Range("H122").Value = 8263.666999999999
$ws.Range("I122").Value = 8869.608
$ws.Range("K122").Value = 26608.824
$ws.Range("M122").Value = -24158.824
$ws.Range("H132").Value = 5412.136
$ws.Range("I132").Value = 5240.421
$ws.Range("J132").Value = 6499.6665
$ws.Range("K132").Value = 15721.263
$ws.Range("L132").Value = 19498.9995
$ws.Range("M132").Value = -13191.263
$ws.Range("N132").Value = -24558.9995
$ws.Range("H136").Value = 7688
$ws.Range("I136").Value = 7027.5713
$ws.Range("K136").Value = 21082.7139
$ws.Range("M136").Value = -18532.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 68997.8
$ws.Range("J109").Value = 68997.8
$ws.Range("L109").Value = 68997.8
$ws.Range("N109").Value = -71771.8
$ws.Range("H126").Value = 558978
$ws.Range("I126").Value = 3828.8572
$ws.Range("K126").Value = 11486.5716
$ws.Range("M126").Value = -9016.571599999999
